$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'47.021.67"
$ws.Range("E2").Value = "  +0.62%  "

# Row 3
$ws.Range("D3").Value = "'2.478.30"
$ws.Range("E3").Value = "  +0.31%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'319.47"
$ws.Range("E5").Value = "  -1.16%  "

# Row 6
$ws.Range("D6").Value = "'107.64"
$ws.Range("E6").Value = "  +2.31%  "

# Row 7
$ws.Range("D7").Value = "'0.520"
$ws.Range("E7").Value = "  -0.43%  "

# Row 8
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
$ws.Range("D9").Value = "'0.531"
$ws.Range("E9").Value = "  -1.25%  "

# Row 10
$ws.Range("D10").Value = "'38.63"
$ws.Range("E10").Value = "  +6.96%  "

# Row 11
$ws.Range("D11").Value = "'0.0806"
$ws.Range("E11").Value = "  -1.14%  "

# Row 12
$ws.Range("E12").Value = "  +0.23%  "

# Row 13
$ws.Range("D13").Value = "'18.07"
$ws.Range("E13").Value = "  -1.10%  "

# Row 14
$ws.Range("D14").Value = "'7.09"
$ws.Range("E14").Value = "  -0.13%  "

# Row 15
$ws.Range("D15").Value = "'2.864.89"
$ws.Range("E15").Value = "  +0.15%  "

# Row 16
$ws.Range("D16").Value = "'2.471.62"
$ws.Range("E16").Value = "  -1.27%  "

# Row 17
$ws.Range("D17").Value = "'0.843"
$ws.Range("E17").Value = "  -0.15%  "

# Row 18
$ws.Range("D18").Value = "'46.890.40"
$ws.Range("E18").Value = "  +0.55%  "

# Row 19
$ws.Range("D19").Value = "'12.69"
$ws.Range("E19").Value = "  +0.14%  "

# Row 20
$ws.Range("D20").Value = "'6.58"
$ws.Range("E20").Value = "  +1.81%  "

# Row 21
$ws.Range("D21").Value = "'2.76"
$ws.Range("E21").Value = "  +15.23%  "

# Row 22
$ws.Range("D22").Value = "'0.0₃0928"
$ws.Range("E22").Value = "  -0.81%  "

# Row 23
$ws.Range("D23").Value = "'70.23"
$ws.Range("E23").Value = "  -0.52%  "

# Row 24
$ws.Range("D24").Value = "'244.49"
$ws.Range("E24").Value = "  -1.88%  "

# Row 25
$ws.Range("D25").Value = "'2.54"
$ws.Range("E25").Value = "  -0.57%  "

# Row 26
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.04%  "

# Row 27
$ws.Range("D27").Value = "'25.52"
$ws.Range("E27").Value = "  -2.47%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.26"
$ws.Range("E28").Value = "  +2.62%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'10.01"
$ws.Range("E29").Value = "  +2.01%  "

# Row 30
$ws.Range("D30").Value = "'0.138"
$ws.Range("E30").Value = "  +4.41%  "

# Row 31
$ws.Range("D31").Value = "'34.77"
$ws.Range("E31").Value = "  -1.21%  "

# Row 32
$ws.Range("D32").Value = "'49.34"
$ws.Range("E32").Value = "  -0.60%  "

# Row 33
$ws.Range("D33").Value = "'19.84"
$ws.Range("E33").Value = "  +1.13%  "

# Row 34
$ws.Range("D34").Value = "'5.31"
$ws.Range("E34").Value = "  -0.08%  "

# Row 35
$ws.Range("D35").Value = "'0.0777"
$ws.Range("E35").Value = "  +1.23%  "

# Row 36
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.14%  "

# Row 37
$ws.Range("D37").Value = "'1.95"
$ws.Range("E37").Value = "  +2.21%  "

# Row 38
$ws.Range("D38").Value = "'4.62"
$ws.Range("E38").Value = "  -0.05%  "

# Row 39
$ws.Range("D39").Value = "'2.93"
$ws.Range("E39").Value = "  -0.65%  "

# Row 40
$ws.Range("E40").Value = "  -0.07%  "

# Row 41
$ws.Range("D41").Value = "'2.22"
$ws.Range("E41").Value = "  -0.22%  "

# Row 42
$ws.Range("D42").Value = "'118.83"
$ws.Range("E42").Value = "  -3.70%  "

# Row 43
$ws.Range("D43").Value = "'21.77"
$ws.Range("E43").Value = "  +5.18%  "

# Row 44
$ws.Range("D44").Value = "'0.0293"
$ws.Range("E44").Value = "  -0.55%  "

# Row 45
$ws.Range("D45").Value = "'1.973.40"
$ws.Range("E45").Value = "  -0.55%  "

# Row 46
$ws.Range("D46").Value = "'2.99"
$ws.Range("E46").Value = "  +0.42%  "

# Row 47
$ws.Range("E47").Value = "  -3.74%  "

# Row 48
$ws.Range("D48").Value = "'9.01"
$ws.Range("E48").Value = "  +0.59%  "

# Row 49
$ws.Range("D49").Value = "'1.75"
$ws.Range("E49").Value = "  -3.16%  "

# Row 50
$ws.Range("D50").Value = "'5.11"
$ws.Range("E50").Value = "  -4.88%  "

# Row 51
$ws.Range("D51").Value = "'56.69"
$ws.Range("E51").Value = "  +3.90%  "
